# Rename the "movies" sheet to "movie", and adjust sheet selection state
# (select C2 on "animated-tv-series", then make "movie" the active sheet)
$wb = $excel.ActiveWorkbook

$moviesSheet = $wb.Worksheets.Item("movies")
$animatedSheet = $wb.Worksheets.Item("animated-tv-series")

# Select a different cell on the animated-tv-series sheet before switching away
$animatedSheet.Activate()
$animatedSheet.Range("C2").Select()

# Rename the sheet
$moviesSheet.Name = "movie"

# Make the renamed sheet the active sheet/tab
$moviesSheet.Activate()
